$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously "latest" row (row 10) is no longer today's row, so it
# switches from the date-only format to the full datetime format.
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new daily row (row 11) with today's data, using the
# date-only format that the last row always carries.
$ws.Range("A11").Value = 45960
$ws.Range("A11").NumberFormat = "YYYY-MM-DD"
$ws.Range("B11").Value = 20
$ws.Range("C11").Value = 29
$ws.Range("D11").Value = 23
